$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tradeDetails")

# tradeType changed from "Temporary" to "Permanent"
$ws.Range("C2").Value = "Permanent"

# tradeCommencementDate changed from "31/03/2017" to "01/03/2017"
$ws.Range("H2").Value = "01/03/2017"

# New row 4 added with tradeCategory/tradeSubCategory data
$ws.Range("D4").Value = "Flammables"
$ws.Range("E4").Value = "Acetylene Gas"
$ws.Range("D4:E4").NumberFormat = "@"

# Update the active selection on the sheet
$ws.Range("H10").Select()
